# Generate Report for Handback
# Updates the timestamp values recorded on the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" on the Overview sheet and the matching
# "Correspond Handoff Datetime" on the de-de sheet share the same value.
$wsOverview.Range("G2").Value = "2016-08-12 07:20:41"
$wsDeDe.Range("H2").Value = "2016-08-12 07:20:41"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-12 07:20:34"
$wsZhCn.Range("K2").Value = "2016-08-12 07:21:05"

# de-de sheet: Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-08-12 07:21:16"
